# ===================================================================
# Applies the KHL stats refresh described in the commit diff:
#  - Matches_SOG: append 3 new match rows (316-318), extend dimension
#  - Shots_HA / Shots_Summary: refresh as_of_utc + recompute stats for
#    the teams involved in the newly-added matches
#  - Meta_ext: bump as_of_utc + build_version
# ===================================================================

$wb = $excel.ActiveWorkbook

# --- Matches_SOG: append new match rows 316-318 ---
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

# Row 316
$wsMatches.Range("A316").NumberFormat = "@"
$wsMatches.Range("A316").Value = "897811"
$wsMatches.Range("A316").Style = "Normal"
$wsMatches.Range("B316").Value = "2025-11-24T19:00:00"
$wsMatches.Range("C316").Value = "Локомотив"
$wsMatches.Range("D316").Value = "Сибирь"
$wsMatches.Range("E316").Value = 23
$wsMatches.Range("F316").Value = 31
$wsMatches.Range("G316").Value = "khl_text"

# Row 317
$wsMatches.Range("A317").NumberFormat = "@"
$wsMatches.Range("A317").Value = "897809"
$wsMatches.Range("A317").Style = "Normal"
$wsMatches.Range("B317").Value = "2025-11-24T19:30:00"
$wsMatches.Range("C317").Value = "Динамо М"
$wsMatches.Range("D317").Value = "Амур"
$wsMatches.Range("E317").Value = 31
$wsMatches.Range("F317").Value = 31
$wsMatches.Range("G317").Value = "khl_text"

# Row 318
$wsMatches.Range("A318").NumberFormat = "@"
$wsMatches.Range("A318").Value = "897810"
$wsMatches.Range("A318").Style = "Normal"
$wsMatches.Range("B318").Value = "2025-11-24T19:30:00"
$wsMatches.Range("C318").Value = "ЦСКА"
$wsMatches.Range("D318").Value = "СКА"
$wsMatches.Range("E318").Value = 23
$wsMatches.Range("F318").Value = 23

# --- Shots_HA: refresh as_of_utc (col D) and recomputed totals ---
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")
$wsShotsHA.Range("D2").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D3").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D4").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D5").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D6").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("F6").Value = 15
$wsShotsHA.Range("K6").Value = 410
$wsShotsHA.Range("L6").Value = 572
$wsShotsHA.Range("M6").Value = 27.3
$wsShotsHA.Range("N6").Value = 38.1
$wsShotsHA.Range("D7").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D8").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("E8").Value = 14
$wsShotsHA.Range("G8").Value = 467
$wsShotsHA.Range("H8").Value = 373
$wsShotsHA.Range("I8").Value = 33.4
$wsShotsHA.Range("J8").Value = 26.6
$wsShotsHA.Range("D9").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D10").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D11").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D12").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("E12").Value = 16
$wsShotsHA.Range("G12").Value = 524
$wsShotsHA.Range("H12").Value = 421
$wsShotsHA.Range("I12").Value = 32.8
$wsShotsHA.Range("J12").Value = 26.3
$wsShotsHA.Range("D13").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D14").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D15").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("F15").Value = 12
$wsShotsHA.Range("K15").Value = 360
$wsShotsHA.Range("L15").Value = 386
$wsShotsHA.Range("M15").Value = 30
$wsShotsHA.Range("N15").Value = 32.2
$wsShotsHA.Range("D16").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D17").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D18").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("F18").Value = 18
$wsShotsHA.Range("K18").Value = 472
$wsShotsHA.Range("L18").Value = 607
$wsShotsHA.Range("M18").Value = 26.2
$wsShotsHA.Range("N18").Value = 33.7
$wsShotsHA.Range("D19").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D20").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D21").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D22").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("D23").Value = "2025-11-24T19:30:00Z"
$wsShotsHA.Range("E23").Value = 13
$wsShotsHA.Range("G23").Value = 312
$wsShotsHA.Range("H23").Value = 367
$wsShotsHA.Range("I23").Value = 24
$wsShotsHA.Range("J23").Value = 28.2

# --- Shots_Summary: refresh as_of_utc (col D) and recomputed totals ---
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")
$wsShotsSummary.Range("D2").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D3").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D4").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D5").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D6").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("E6").Value = 29
$wsShotsSummary.Range("F6").Value = 829
$wsShotsSummary.Range("G6").Value = 1062
$wsShotsSummary.Range("H6").Value = 28.6
$wsShotsSummary.Range("I6").Value = 36.6
$wsShotsSummary.Range("D7").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D8").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("E8").Value = 28
$wsShotsSummary.Range("F8").Value = 852
$wsShotsSummary.Range("G8").Value = 809
$wsShotsSummary.Range("I8").Value = 28.9
$wsShotsSummary.Range("D9").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D10").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D11").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D12").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("E12").Value = 31
$wsShotsSummary.Range("F12").Value = 993
$wsShotsSummary.Range("G12").Value = 791
$wsShotsSummary.Range("H12").Value = 32
$wsShotsSummary.Range("I12").Value = 25.5
$wsShotsSummary.Range("D13").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D14").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D15").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("E15").Value = 27
$wsShotsSummary.Range("F15").Value = 848
$wsShotsSummary.Range("G15").Value = 882
$wsShotsSummary.Range("H15").Value = 31.4
$wsShotsSummary.Range("I15").Value = 32.7
$wsShotsSummary.Range("D16").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D17").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D18").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("E18").Value = 30
$wsShotsSummary.Range("F18").Value = 800
$wsShotsSummary.Range("G18").Value = 1041
$wsShotsSummary.Range("H18").Value = 26.7
$wsShotsSummary.Range("I18").Value = 34.7
$wsShotsSummary.Range("D19").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D20").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D21").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D22").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("D23").Value = "2025-11-24T19:30:00Z"
$wsShotsSummary.Range("E23").Value = 29
$wsShotsSummary.Range("F23").Value = 726
$wsShotsSummary.Range("G23").Value = 792
$wsShotsSummary.Range("H23").Value = 25
$wsShotsSummary.Range("I23").Value = 27.3

# --- Meta_ext: bump as_of_utc and build_version ---
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-24T19:30:00Z"
$wsMeta.Range("D2").Value = 11

Write-Host "Applied KHL stats refresh: Matches_SOG +3 rows, Shots_HA/Shots_Summary/Meta_ext updated."